$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D:E columns for data rows remain text so numeric-looking strings are preserved
$ws.Range("D2:E51").NumberFormat = "@"

$sub3 = [string]([char]0x2083)

$ws.Range("D2").Value = "36.102.96"
$ws.Range("E2").Value = "  -4.08%  "

$ws.Range("D3").Value = "1.965.22"
$ws.Range("E3").Value = "  -3.33%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "228.82"
$ws.Range("E5").Value = "  -14.24%  "

$ws.Range("D6").Value = "0.594"
$ws.Range("E6").Value = "  -4.51%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "53.47"
$ws.Range("E8").Value = "  -6.79%  "

$ws.Range("D9").Value = "0.365"
$ws.Range("E9").Value = "  -6.01%  "

$ws.Range("D10").Value = "57.66"
$ws.Range("E10").Value = "  +0.97%  "

$ws.Range("D11").Value = "0.0739"
$ws.Range("E11").Value = "  -6.10%  "

$ws.Range("D12").Value = "0.0972"
$ws.Range("E12").Value = "  -4.91%  "

$ws.Range("D13").Value = "2.256.31"
$ws.Range("E13").Value = "  -3.18%  "

$ws.Range("D14").Value = "13.69"
$ws.Range("E14").Value = "  -6.30%  "

$ws.Range("D15").Value = "19.74"
$ws.Range("E15").Value = "  -6.35%  "

$ws.Range("D16").Value = "0.743"
$ws.Range("E16").Value = "  -9.25%  "

$ws.Range("D17").Value = "4.97"
$ws.Range("E17").Value = "  -6.78%  "

$ws.Range("D18").Value = "1.969.97"
$ws.Range("E18").Value = "  -3.49%  "

$ws.Range("D19").Value = "36.158.94"
$ws.Range("E19").Value = "  -3.67%  "

$ws.Range("D20").Value = "67.24"
$ws.Range("E20").Value = "  -4.37%  "

$ws.Range("D21").Value = "0.0" + $sub3 + "0796"
$ws.Range("E21").Value = "  -6.61%  "

$ws.Range("D22").Value = "4.99"
$ws.Range("E22").Value = "  -4.79%  "

$ws.Range("D23").Value = "220.08"
$ws.Range("E23").Value = "  -3.67%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  -16.39%  "

$ws.Range("D27").Value = "158.86"
$ws.Range("E27").Value = "  -3.61%  "

$ws.Range("D28").Value = "8.45"
$ws.Range("E28").Value = "  -7.06%  "

$ws.Range("D29").Value = "18.59"
$ws.Range("E29").Value = "  -6.32%  "

$ws.Range("D30").Value = "1.31"
$ws.Range("E30").Value = "  -3.71%  "

$ws.Range("E31").Value = "  -6.80%  "

$ws.Range("E32").Value = "  -4.32%  "

$ws.Range("E33").Value = "  -8.43%  "

$ws.Range("D34").Value = "0.0594"
$ws.Range("E34").Value = "  -10.59%  "

$ws.Range("D35").Value = "4.19"
$ws.Range("E35").Value = "  -8.48%  "

$ws.Range("D36").Value = "2.26"
$ws.Range("E36").Value = "  -7.50%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("D39").Value = "3.16"
$ws.Range("E39").Value = "  -6.85%  "

$ws.Range("E40").Value = "  -3.20%  "

$ws.Range("D41").Value = "3.00"
$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("D42").Value = "1.416.41"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").Value = "0.0875"
$ws.Range("E43").Value = "  -8.34%  "

$ws.Range("D44").Value = "0.0198"
$ws.Range("E44").Value = "  -8.46%  "

$ws.Range("E45").Value = "  -13.49%  "

$ws.Range("D46").Value = "86.56"
$ws.Range("E46").Value = "  -5.87%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.974"
$ws.Range("E47").Value = "  -7.43%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "2.86"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "14.49"
$ws.Range("E49").Value = "  -9.50%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "6.66"
$ws.Range("E50").Value = "  -7.43%  "

$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "3.55"
$ws.Range("E51").Value = "  +10.26%  "
